$wb = $excel.ActiveWorkbook

# 1. Duplicate the "EJ45" sheet, placing the copy immediately after it.
#    Excel names the duplicate "EJ45 (2)" automatically.
$ej45 = $wb.Worksheets.Item("EJ45")
[void]$ej45.Copy($null, $ej45)
$ej45copy = $wb.Worksheets.Item("EJ45 (2)")

# 2. Update the measurement data on the new copy (columns B, D, G only;
#    A, C, E, F stay identical to the original EJ45 sheet).
$ej45copy.Range("B2").Value = 328.6
$ej45copy.Range("D2").Value = 191.6
$ej45copy.Range("G2").Value = 80.9

$ej45copy.Range("B3").Value = 334.3
$ej45copy.Range("D3").Value = 199.1
$ej45copy.Range("G3").Value = 83.2

$ej45copy.Range("B4").Value = 340.7
$ej45copy.Range("D4").Value = 205.7
$ej45copy.Range("G4").Value = 85.5

$ej45copy.Range("B5").Value = 346.7
$ej45copy.Range("D5").Value = 212.7
$ej45copy.Range("G5").Value = 88.2

$ej45copy.Range("B6").Value = 352.5
$ej45copy.Range("D6").Value = 216
$ej45copy.Range("G6").Value = 90.9

$ej45copy.Range("B7").Value = 359.2
$ej45copy.Range("D7").Value = 220.5
$ej45copy.Range("G7").Value = 93.8

$ej45copy.Range("B8").Value = 366.7
$ej45copy.Range("D8").Value = 227.3
$ej45copy.Range("G8").Value = 95.9

# Leave the cursor where typing the last value would have left it.
[void]$ej45copy.Range("G9").Select()

# 3. Rename the header row of the original "EJ45" sheet to the
#    string-specific labels, and make it the active sheet/cell again.
$ej45.Range("B1").Value = "J4501"
$ej45.Range("C1").Value = "J4502"
$ej45.Range("D1").Value = "J4503"
$ej45.Range("E1").Value = "J4504"
$ej45.Range("F1").Value = "J4505"
$ej45.Range("G1").Value = "J4506"

[void]$ej45.Select()
[void]$ej45.Range("B1").Select()
